# Auto-generated Excel COM-interop script to apply market-price/profit updates
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3300.44
$ws.Range("I40").Value = 4157.5
$ws.Range("J40").Value = 2897.1177
$ws.Range("K40").Value = 4157.5
$ws.Range("L40").Value = 2897.1177
$ws.Range("M40").Value = -3982.5
$ws.Range("N40").Value = -3247.1177
$ws.Range("H82").Value = 3762.3157
$ws.Range("I82").Value = 1074.8572
$ws.Range("J82").Value = 5330
$ws.Range("K82").Value = 3224.5716
$ws.Range("L82").Value = 15990
$ws.Range("M82").Value = -2818.5716
$ws.Range("N82").Value = -16802
$ws.Range("H85").Value = 3762.3157
$ws.Range("I85").Value = 1074.8572
$ws.Range("J85").Value = 5330
$ws.Range("K85").Value = 3224.5716
$ws.Range("L85").Value = 15990
$ws.Range("M85").Value = -1820.5716
$ws.Range("N85").Value = -18798
$ws.Range("H88").Value = 3349.0715
$ws.Range("I88").Value = 1720.75
$ws.Range("J88").Value = 4000.4
$ws.Range("K88").Value = 1720.75
$ws.Range("L88").Value = 4000.4
$ws.Range("M88").Value = -1314.75
$ws.Range("N88").Value = -4812.4
$ws.Range("H91").Value = 3349.0715
$ws.Range("I91").Value = 1720.75
$ws.Range("J91").Value = 4000.4
$ws.Range("K91").Value = 1720.75
$ws.Range("L91").Value = 4000.4
$ws.Range("M91").Value = -316.75
$ws.Range("N91").Value = -6808.4
$ws.Range("H132").Value = 3972944.2
$ws.Range("I132").Value = 4241842
$ws.Range("K132").Value = 12725526
$ws.Range("M132").Value = -12722996
$ws.Range("H135").Value = 1197.9423
$ws.Range("I135").Value = 470.69232
$ws.Range("J135").Value = 3379.6924
$ws.Range("K135").Value = 4236.23088
$ws.Range("L135").Value = 30417.2316
$ws.Range("M135").Value = -1701.23088
$ws.Range("N135").Value = -35487.2316
$ws.Range("H138").Value = 2178.92
$ws.Range("I138").Value = 1344.9032
$ws.Range("J138").Value = 2553.6233
$ws.Range("K138").Value = 4034.7096
$ws.Range("L138").Value = 7660.869900000001
$ws.Range("M138").Value = 1105.2904
$ws.Range("N138").Value = -17940.8699
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 84803.586
$ws.Range("I45").Value = 167633.17
$ws.Range("J45").Value = 1974
$ws.Range("K45").Value = 167633.17
$ws.Range("L45").Value = 1974
$ws.Range("M45").Value = -167256.17
$ws.Range("N45").Value = -2728
$ws.Range("H61").Value = 2121.3794
$ws.Range("I61").Value = 2095.5454
$ws.Range("J61").Value = 2137.1667
$ws.Range("K61").Value = 2095.5454
$ws.Range("L61").Value = 2137.1667
$ws.Range("M61").Value = -1883.5454
$ws.Range("N61").Value = -2561.1667
$ws.Range("H63").Value = 2396.6667
$ws.Range("I63").Value = 1200
$ws.Range("J63").Value = 2995
$ws.Range("K63").Value = 1200
$ws.Range("L63").Value = 2995
$ws.Range("M63").Value = -514
$ws.Range("N63").Value = -4367
$ws.Range("H66").Value = 2396.6667
$ws.Range("I66").Value = 1200
$ws.Range("J66").Value = 2995
$ws.Range("K66").Value = 6000
$ws.Range("L66").Value = 14975
$ws.Range("M66").Value = -2568
$ws.Range("N66").Value = -21839
$ws.Range("H74").Value = 1135.579
$ws.Range("I74").Value = 980.7692
$ws.Range("K74").Value = 980.7692
$ws.Range("M74").Value = -106.7692
$ws.Range("H77").Value = 1135.579
$ws.Range("I77").Value = 980.7692
$ws.Range("K77").Value = 4903.846
$ws.Range("M77").Value = -535.8459999999995
$ws.Range("H132").Value = 9668.338
$ws.Range("I132").Value = 12313.412
$ws.Range("K132").Value = 36940.236
$ws.Range("M132").Value = -34410.236
$ws.Range("H136").Value = 2121.3794
$ws.Range("I136").Value = 2095.5454
$ws.Range("J136").Value = 2137.1667
$ws.Range("K136").Value = 6286.6362
$ws.Range("L136").Value = 6411.500100000001
$ws.Range("M136").Value = -3736.6362
$ws.Range("N136").Value = -11511.5001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 3500
$ws.Range("J15").Value = 3500
$ws.Range("L15").Value = 3500
$ws.Range("N15").Value = -3954
$ws.Range("H82").Value = 18122.5
$ws.Range("J82").Value = 33081
$ws.Range("L82").Value = 33081
$ws.Range("N82").Value = -33847
$ws.Range("H85").Value = 18122.5
$ws.Range("J85").Value = 33081
$ws.Range("L85").Value = 33081
$ws.Range("N85").Value = -35733
$ws.Range("H99").Value = 1918.6875
$ws.Range("I99").Value = 1931.6666
$ws.Range("K99").Value = 1931.6666
$ws.Range("M99").Value = -433.6666
$ws.Range("H134").Value = 3263.204
$ws.Range("I134").Value = 3579.9375
$ws.Range("J134").Value = 2667
$ws.Range("K134").Value = 10739.8125
$ws.Range("L134").Value = 8001
$ws.Range("M134").Value = -8204.8125
$ws.Range("N134").Value = -13071
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28026.982
$ws.Range("I31").Value = 1467.0555
$ws.Range("J31").Value = 40285.41
$ws.Range("K31").Value = 1467.0555
$ws.Range("L31").Value = 40285.41
$ws.Range("M31").Value = -1172.0555
$ws.Range("N31").Value = -40875.41
$ws.Range("H34").Value = 28026.982
$ws.Range("I34").Value = 1467.0555
$ws.Range("J34").Value = 40285.41
$ws.Range("K34").Value = 1467.0555
$ws.Range("L34").Value = 40285.41
$ws.Range("M34").Value = -1265.0555
$ws.Range("N34").Value = -40689.41
$ws.Range("H58").Value = 1669.4359
$ws.Range("I58").Value = 1454
$ws.Range("J58").Value = 2100.3076
$ws.Range("K58").Value = 1454
$ws.Range("L58").Value = 2100.3076
$ws.Range("M58").Value = -1251
$ws.Range("N58").Value = -2506.3076
$ws.Range("H59").Value = 26940
$ws.Range("J59").Value = 26940
$ws.Range("L59").Value = 26940
$ws.Range("N59").Value = -29230
$ws.Range("H68").Value = 17623.705
$ws.Range("J68").Value = 17623.705
$ws.Range("L68").Value = 17623.705
$ws.Range("N68").Value = -19121.705
$ws.Range("H71").Value = 17623.705
$ws.Range("J71").Value = 17623.705
$ws.Range("L71").Value = 52871.11500000001
$ws.Range("N71").Value = -60359.11500000001
$ws.Range("H74").Value = 39689.75
$ws.Range("J74").Value = 39689.75
$ws.Range("L74").Value = 39689.75
$ws.Range("N74").Value = -41437.75
$ws.Range("H77").Value = 39689.75
$ws.Range("J77").Value = 39689.75
$ws.Range("L77").Value = 119069.25
$ws.Range("N77").Value = -127805.25
$ws.Range("H132").Value = 2939.2415
$ws.Range("I132").Value = 2818.7273
$ws.Range("J132").Value = 3318
$ws.Range("K132").Value = 8456.1819
$ws.Range("L132").Value = 9954
$ws.Range("M132").Value = -5926.1819
$ws.Range("N132").Value = -15014
$ws.Range("H134").Value = 1062.6875
$ws.Range("I134").Value = 642.4167
$ws.Range("J134").Value = 2323.5
$ws.Range("K134").Value = 1927.2501
$ws.Range("L134").Value = 6970.5
$ws.Range("M134").Value = 607.7499
$ws.Range("N134").Value = -12040.5
$ws.Range("H136").Value = 1669.4359
$ws.Range("I136").Value = 1454
$ws.Range("J136").Value = 2100.3076
$ws.Range("K136").Value = 4362
$ws.Range("L136").Value = 6300.9228
$ws.Range("M136").Value = -1812
$ws.Range("N136").Value = -11400.9228
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 50.833332
$ws.Range("J12").Value = 54.3125
$ws.Range("L12").Value = 162.9375
$ws.Range("N12").Value = -508.9375
$ws.Range("H22").Value = 3099.9805
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 900
$ws.Range("M22").Value = -731
$ws.Range("H27").Value = 3099.9805
$ws.Range("I27").Value = 300
$ws.Range("K27").Value = 900
$ws.Range("M27").Value = -798
$ws.Range("H39").Value = 5200
$ws.Range("J39").Value = 6766.6665
$ws.Range("L39").Value = 20299.9995
$ws.Range("N39").Value = -20887.9995
$ws.Range("H58").Value = 2510.2856
$ws.Range("I58").Value = 1077.5
$ws.Range("J58").Value = 2847.4119
$ws.Range("K58").Value = 3232.5
$ws.Range("L58").Value = 8542.235700000001
$ws.Range("M58").Value = -3104.5
$ws.Range("N58").Value = -8798.235700000001
$ws.Range("H86").Value = 1116.6666
$ws.Range("I86").Value = 1266.6666
$ws.Range("J86").Value = 966.6667
$ws.Range("K86").Value = 3799.9998
$ws.Range("L86").Value = 2900.0001
$ws.Range("M86").Value = -2613.9998
$ws.Range("N86").Value = -5272.0001
$ws.Range("H89").Value = 1116.6666
$ws.Range("I89").Value = 1266.6666
$ws.Range("J89").Value = 966.6667
$ws.Range("K89").Value = 11399.9994
$ws.Range("L89").Value = 8700.0003
$ws.Range("M89").Value = -5471.999400000001
$ws.Range("N89").Value = -20556.0003
$ws.Range("H92").Value = 1000.75
$ws.Range("I92").Value = 300
$ws.Range("J92").Value = 1234.3334
$ws.Range("K92").Value = 900
$ws.Range("L92").Value = 3703.0002
$ws.Range("M92").Value = 348
$ws.Range("N92").Value = -6199.0002
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2468.451
$ws.Range("I132").Value = 1987.2812
$ws.Range("J132").Value = 3278.842
$ws.Range("K132").Value = 5961.8436
$ws.Range("L132").Value = 9836.526
$ws.Range("M132").Value = -3431.8436
$ws.Range("N132").Value = -14896.526
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4827.0713
$ws.Range("I46").Value = 4234.875
$ws.Range("J46").Value = 5616.6665
$ws.Range("K46").Value = 4234.875
$ws.Range("L46").Value = 5616.6665
$ws.Range("M46").Value = -4046.875
$ws.Range("N46").Value = -5992.6665
$ws.Range("H132").Value = 3636.1892
$ws.Range("I132").Value = 4326.8096
$ws.Range("J132").Value = 2729.75
$ws.Range("K132").Value = 12980.4288
$ws.Range("L132").Value = 8189.25
$ws.Range("M132").Value = -10450.4288
$ws.Range("N132").Value = -13249.25
$ws.Range("H136").Value = 2298.5715
$ws.Range("I136").Value = 1517.2174
$ws.Range("J136").Value = 5892.8
$ws.Range("K136").Value = 4551.6522
$ws.Range("L136").Value = 17678.4
$ws.Range("M136").Value = -2001.6522
$ws.Range("N136").Value = -22778.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 618.35297
$ws.Range("I113").Value = 489.3846
$ws.Range("K113").Value = 1468.1538
$ws.Range("M113").Value = 701.8462
$ws.Range("H136").Value = 15814.682
$ws.Range("I136").Value = 28288.445
$ws.Range("K136").Value = 84865.33499999999
$ws.Range("M136").Value = -82315.33499999999
